$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-06-26 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-27 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("55-53=2", $true, $false, $false, $false, $false, $true, 1, $false, "55+17=72", 2) | Out-Null
$d.Content.Find.Execute("36+42=78", $true, $false, $false, $false, $false, $true, 1, $false, "47-32=15", 2) | Out-Null
$d.Content.Find.Execute("9+52=61", $true, $false, $false, $false, $false, $true, 1, $false, "22-21=1", 2) | Out-Null
$d.Content.Find.Execute("49+13=62", $true, $false, $false, $false, $false, $true, 1, $false, "97-68=29", 2) | Out-Null
$d.Content.Find.Execute("48-7=41", $true, $false, $false, $false, $false, $true, 1, $false, "82-55=27", 2) | Out-Null
$d.Content.Find.Execute("25+14=39", $true, $false, $false, $false, $false, $true, 1, $false, "54+38=92", 2) | Out-Null
$d.Content.Find.Execute("14-4=10", $true, $false, $false, $false, $false, $true, 1, $false, "89+1=90", 2) | Out-Null
$d.Content.Find.Execute("61+29=90", $true, $false, $false, $false, $false, $true, 1, $false, "51-4=47", 2) | Out-Null
$d.Content.Find.Execute("45+27=72", $true, $false, $false, $false, $false, $true, 1, $false, "72+8=80", 2) | Out-Null
$d.Content.Find.Execute("87-46=41", $true, $false, $false, $false, $false, $true, 1, $false, "75-10=65", 2) | Out-Null
$d.Content.Find.Execute("74-28=46", $true, $false, $false, $false, $false, $true, 1, $false, "85+9=94", 2) | Out-Null
$d.Content.Find.Execute("55-26=29", $true, $false, $false, $false, $false, $true, 1, $false, "66+30=96", 2) | Out-Null
$d.Content.Find.Execute("22-18=4", $true, $false, $false, $false, $false, $true, 1, $false, "93-5=88", 2) | Out-Null
$d.Content.Find.Execute("50+25=75", $true, $false, $false, $false, $false, $true, 1, $false, "79-25=54", 2) | Out-Null
$d.Content.Find.Execute("99-83=16", $true, $false, $false, $false, $false, $true, 1, $false, "48-3=45", 2) | Out-Null
$d.Content.Find.Execute("51+3=54", $true, $false, $false, $false, $false, $true, 1, $false, "40+2=42", 2) | Out-Null
$d.Content.Find.Execute("77-46=31", $true, $false, $false, $false, $false, $true, 1, $false, "70-19=51", 2) | Out-Null
$d.Content.Find.Execute("68+15=83", $true, $false, $false, $false, $false, $true, 1, $false, "99-1=98", 2) | Out-Null
$d.Content.Find.Execute("8+67=75", $true, $false, $false, $false, $false, $true, 1, $false, "60-43=17", 2) | Out-Null
$d.Content.Find.Execute("59-6=53", $true, $false, $false, $false, $false, $true, 1, $false, "71-29=42", 2) | Out-Null
$d.Content.Find.Execute("79-66=13", $true, $false, $false, $false, $false, $true, 1, $false, "79-12=67", 2) | Out-Null
$d.Content.Find.Execute("21+62=83", $true, $false, $false, $false, $false, $true, 1, $false, "95-81=14", 2) | Out-Null
$d.Content.Find.Execute("25+47=72", $true, $false, $false, $false, $false, $true, 1, $false, "78-69=9", 2) | Out-Null
$d.Content.Find.Execute("24+38=62", $true, $false, $false, $false, $false, $true, 1, $false, "41+23=64", 2) | Out-Null
$d.Content.Find.Execute("77-73=4", $true, $false, $false, $false, $false, $true, 1, $false, "99-45=54", 2) | Out-Null
$d.Content.Find.Execute("44+14=58", $true, $false, $false, $false, $false, $true, 1, $false, "52-18=34", 2) | Out-Null
$d.Content.Find.Execute("34+9=43", $true, $false, $false, $false, $false, $true, 1, $false, "78-63=15", 2) | Out-Null
$d.Content.Find.Execute("97-42=55", $true, $false, $false, $false, $false, $true, 1, $false, "54-36=18", 2) | Out-Null
$d.Content.Find.Execute("24+42=66", $true, $false, $false, $false, $false, $true, 1, $false, "51-30=21", 2) | Out-Null
$d.Content.Find.Execute("2+36=38", $true, $false, $false, $false, $false, $true, 1, $false, "31+41=72", 2) | Out-Null
$d.Content.Find.Execute("45+23=68", $true, $false, $false, $false, $false, $true, 1, $false, "24+7=31", 2) | Out-Null
$d.Content.Find.Execute("91-66=25", $true, $false, $false, $false, $false, $true, 1, $false, "34+41=75", 2) | Out-Null
$d.Content.Find.Execute("64-11=53", $true, $false, $false, $false, $false, $true, 1, $false, "24+51=75", 2) | Out-Null
$d.Content.Find.Execute("67-14=53", $true, $false, $false, $false, $false, $true, 1, $false, "34+11=45", 2) | Out-Null
$d.Content.Find.Execute("69-17=52", $true, $false, $false, $false, $false, $true, 1, $false, "0+2=2", 2) | Out-Null
$d.Content.Find.Execute("77-27=50", $true, $false, $false, $false, $false, $true, 1, $false, "79-27=52", 2) | Out-Null
$d.Content.Find.Execute("44+51=95", $true, $false, $false, $false, $false, $true, 1, $false, "0+35=35", 2) | Out-Null
$d.Content.Find.Execute("0+12=12", $true, $false, $false, $false, $false, $true, 1, $false, "26-22=4", 2) | Out-Null
$d.Content.Find.Execute("75+21=96", $true, $false, $false, $false, $false, $true, 1, $false, "75+12=87", 2) | Out-Null
$d.Content.Find.Execute("35+17=52", $true, $false, $false, $false, $false, $true, 1, $false, "4+27=31", 2) | Out-Null
$d.Content.Find.Execute("47+47=94", $true, $false, $false, $false, $false, $true, 1, $false, "84-56=28", 2) | Out-Null
$d.Content.Find.Execute("87-34=53", $true, $false, $false, $false, $false, $true, 1, $false, "68+9=77", 2) | Out-Null
$d.Content.Find.Execute("67+6=73", $true, $false, $false, $false, $false, $true, 1, $false, "25+24=49", 2) | Out-Null
$d.Content.Find.Execute("23-20=3", $true, $false, $false, $false, $false, $true, 1, $false, "34+54=88", 2) | Out-Null
$d.Content.Find.Execute("7+49=56", $true, $false, $false, $false, $false, $true, 1, $false, "37+60=97", 2) | Out-Null
$d.Content.Find.Execute("63+19=82", $true, $false, $false, $false, $false, $true, 1, $false, "43+13=56", 2) | Out-Null
$d.Content.Find.Execute("10+68=78", $true, $false, $false, $false, $false, $true, 1, $false, "81-35=46", 2) | Out-Null
$d.Content.Find.Execute("89+2=91", $true, $false, $false, $false, $false, $true, 1, $false, "0+61=61", 2) | Out-Null
$d.Content.Find.Execute("25+25=50", $true, $false, $false, $false, $false, $true, 1, $false, "96-4=92", 2) | Out-Null
$d.Content.Find.Execute("73-47=26", $true, $false, $false, $false, $false, $true, 1, $false, "57-24=33", 2) | Out-Null
$d.Content.Find.Execute("68-45=23", $true, $false, $false, $false, $false, $true, 1, $false, "47+34=81", 2) | Out-Null
$d.Content.Find.Execute("10+33=43", $true, $false, $false, $false, $false, $true, 1, $false, "44+2=46", 2) | Out-Null
$d.Content.Find.Execute("90-52=38", $true, $false, $false, $false, $false, $true, 1, $false, "24+19=43", 2) | Out-Null
$d.Content.Find.Execute("28+42=70", $true, $false, $false, $false, $false, $true, 1, $false, "39+8=47", 2) | Out-Null
$d.Content.Find.Execute("21-16=5", $true, $false, $false, $false, $false, $true, 1, $false, "65+24=89", 2) | Out-Null
$d.Content.Find.Execute("6+88=94", $true, $false, $false, $false, $false, $true, 1, $false, "49-11=38", 2) | Out-Null
$d.Content.Find.Execute("45+29=74", $true, $false, $false, $false, $false, $true, 1, $false, "45+7=52", 2) | Out-Null
$d.Content.Find.Execute("22+50=72", $true, $false, $false, $false, $false, $true, 1, $false, "43+26=69", 2) | Out-Null
$d.Content.Find.Execute("83-79=4", $true, $false, $false, $false, $false, $true, 1, $false, "77+16=93", 2) | Out-Null
$d.Content.Find.Execute("1+21=22", $true, $false, $false, $false, $false, $true, 1, $false, "73-51=22", 2) | Out-Null
$d.Content.Find.Execute("49+26=75", $true, $false, $false, $false, $false, $true, 1, $false, "3+1=4", 2) | Out-Null
$d.Content.Find.Execute("88-49=39", $true, $false, $false, $false, $false, $true, 1, $false, "75+4=79", 2) | Out-Null
$d.Content.Find.Execute("20+17=37", $true, $false, $false, $false, $false, $true, 1, $false, "71-61=10", 2) | Out-Null
$d.Content.Find.Execute("25+66=91", $true, $false, $false, $false, $false, $true, 1, $false, "17+17=34", 2) | Out-Null
$d.Content.Find.Execute("68-36=32", $true, $false, $false, $false, $false, $true, 1, $false, "68+19=87", 2) | Out-Null
$d.Content.Find.Execute("79-5=74", $true, $false, $false, $false, $false, $true, 1, $false, "57-32=25", 2) | Out-Null
$d.Content.Find.Execute("27+31=58", $true, $false, $false, $false, $false, $true, 1, $false, "2+72=74", 2) | Out-Null
$d.Content.Find.Execute("1+90=91", $true, $false, $false, $false, $false, $true, 1, $false, "29+9=38", 2) | Out-Null
$d.Content.Find.Execute("65-47=18", $true, $false, $false, $false, $false, $true, 1, $false, "90+4=94", 2) | Out-Null
$d.Content.Find.Execute("66-3=63", $true, $false, $false, $false, $false, $true, 1, $false, "88-27=61", 2) | Out-Null
$d.Content.Find.Execute("37+57=94", $true, $false, $false, $false, $false, $true, 1, $false, "75-55=20", 2) | Out-Null
$d.Content.Find.Execute("4+41=45", $true, $false, $false, $false, $false, $true, 1, $false, "22-9=13", 2) | Out-Null
$d.Content.Find.Execute("49+21=70", $true, $false, $false, $false, $false, $true, 1, $false, "23+44=67", 2) | Out-Null
$d.Content.Find.Execute("87-38=49", $true, $false, $false, $false, $false, $true, 1, $false, "69-20=49", 2) | Out-Null
$d.Content.Find.Execute("90+1=91", $true, $false, $false, $false, $false, $true, 1, $false, "28-22=6", 2) | Out-Null
$d.Content.Find.Execute("8+59=67", $true, $false, $false, $false, $false, $true, 1, $false, "32+45=77", 2) | Out-Null
$d.Content.Find.Execute("43+17=60", $true, $false, $false, $false, $false, $true, 1, $false, "88-63=25", 2) | Out-Null
$d.Content.Find.Execute("66-41=25", $true, $false, $false, $false, $false, $true, 1, $false, "26+64=90", 2) | Out-Null
$d.Content.Find.Execute("63-58=5", $true, $false, $false, $false, $false, $true, 1, $false, "57-50=7", 2) | Out-Null
$d.Content.Find.Execute("13+62=75", $true, $false, $false, $false, $false, $true, 1, $false, "37+25=62", 2) | Out-Null
$d.Content.Find.Execute("16+44=60", $true, $false, $false, $false, $false, $true, 1, $false, "48-37=11", 2) | Out-Null
$d.Content.Find.Execute("49-29=20", $true, $false, $false, $false, $false, $true, 1, $false, "48-39=9", 2) | Out-Null
$d.Content.Find.Execute("1+43=44", $true, $false, $false, $false, $false, $true, 1, $false, "93-1=92", 2) | Out-Null
$d.Content.Find.Execute("16-7=9", $true, $false, $false, $false, $false, $true, 1, $false, "85-25=60", 2) | Out-Null
$d.Content.Find.Execute("72+5=77", $true, $false, $false, $false, $false, $true, 1, $false, "33+11=44", 2) | Out-Null
$d.Content.Find.Execute("54+22=76", $true, $false, $false, $false, $false, $true, 1, $false, "0+60=60", 2) | Out-Null
$d.Content.Find.Execute("14+9=23", $true, $false, $false, $false, $false, $true, 1, $false, "55-22=33", 2) | Out-Null
$d.Content.Find.Execute("28+65=93", $true, $false, $false, $false, $false, $true, 1, $false, "98-12=86", 2) | Out-Null
$d.Content.Find.Execute("45-2=43", $true, $false, $false, $false, $false, $true, 1, $false, "29-10=19", 2) | Out-Null
$d.Content.Find.Execute("25-17=8", $true, $false, $false, $false, $false, $true, 1, $false, "94-88=6", 2) | Out-Null
$d.Content.Find.Execute("69+30=99", $true, $false, $false, $false, $false, $true, 1, $false, "6+61=67", 2) | Out-Null
$d.Content.Find.Execute("83-13=70", $true, $false, $false, $false, $false, $true, 1, $false, "4+54=58", 2) | Out-Null
$d.Content.Find.Execute("65+21=86", $true, $false, $false, $false, $false, $true, 1, $false, "8+55=63", 2) | Out-Null
$d.Content.Find.Execute("52-50=2", $true, $false, $false, $false, $false, $true, 1, $false, "35+6=41", 2) | Out-Null
$d.Content.Find.Execute("32-12=20", $true, $false, $false, $false, $false, $true, 1, $false, "25+60=85", 2) | Out-Null
$d.Content.Find.Execute("34+39=73", $true, $false, $false, $false, $false, $true, 1, $false, "72-35=37", 2) | Out-Null
$d.Content.Find.Execute("30+57=87", $true, $false, $false, $false, $false, $true, 1, $false, "53+23=76", 2) | Out-Null
$d.Content.Find.Execute("48-46=2", $true, $false, $false, $false, $false, $true, 1, $false, "2+66=68", 2) | Out-Null
$d.Content.Find.Execute("59+40=99", $true, $false, $false, $false, $false, $true, 1, $false, "68-50=18", 2) | Out-Null
$d.Content.Find.Execute("24+32=56", $true, $false, $false, $false, $false, $true, 1, $false, "3+41=44", 2) | Out-Null
